$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Step 1: turn "...bra repetida asi como documento." into
#         "...bra repetida, asi como documento." while forcing the run to
#         split exactly at the word boundaries (matching how Word splits
#         runs around an edit). We do this by temporarily bookmarking the
#         two split points (this is what makes the engine keep the runs
#         separate), inserting the comma, and then removing the helper
#         bookmarks again.
# --------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("bra repetida")
$posAfterBra = $idx + 4                 # right after "bra "
$posAfterRepetida = $posAfterBra + 8    # right after "repetida"

$d.Bookmarks.Add("TmpSplit1", $d.Range($posAfterBra, $posAfterBra))
$d.Bookmarks.Add("TmpSplit2", $d.Range($posAfterRepetida, $posAfterRepetida))

$d.Range($posAfterRepetida, $posAfterRepetida).InsertAfter(",")

$d.Bookmarks("TmpSplit1").Delete()
$d.Bookmarks("TmpSplit2").Delete()

# --------------------------------------------------------------------------
# Step 2: append the new closing sentence at the very end of the first
# paragraph (right before the paragraph mark), as two separate runs so we
# can drop the relocated _GoBack bookmark between them.
# --------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$endPos = $p1.End - 1   # position right before the paragraph mark

$part1 = " Hola mi nombre"
$part2 = " como te decía mi nombre es"

$d.Range($endPos, $endPos).InsertAfter($part1)
$midPos = $endPos + $part1.Length
$d.Range($midPos, $midPos).InsertAfter($part2)

# --------------------------------------------------------------------------
# Step 3: relocate the "_GoBack" bookmark from its original spot to the
# newly-created gap between the two appended runs.
# --------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($midPos, $midPos))

Write-Host "Edit complete"
